# Update the "Generate Date" / handoff-handback timestamps shown in the
# handback status report, as part of generating a fresh report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for the first row.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-03 19:54:28"

# --- zh-cn sheet ------------------------------------------------------
# "Correspond Handoff Datetime" and "Correspond Handback DateTime".
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-03 19:54:14"
$wsZhCn.Range("K2").Value = "2016-11-03 19:55:04"

# --- de-de sheet ------------------------------------------------------
# "Correspond Handoff Datetime" and "Correspond Handback DateTime".
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-03 19:54:28"
$wsDeDe.Range("K2").Value = "2016-11-03 19:55:22"

Write-Output "Updated handback status timestamps."
